# Add "Wins", "Losses", "Ties" season-record columns (AD:AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1 - reuse the exact formatting of the existing header row
# (A1:AC1) by copying A1's format onto the new header cells, then set text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-38).
$wins = 70
$losses = 43
$ties = 0

for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
